$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 398.83334
$ws.Range("I2").Value = 97
$ws.Range("K2").Value = 97
$ws.Range("M2").Value = 16

$ws.Range("H33").Value = 616.0625
$ws.Range("J33").Value = 1048.1666
$ws.Range("L33").Value = 1048.1666
$ws.Range("N33").Value = -1506.1666

$ws.Range("H40").Value = 3847.6191
$ws.Range("I40").Value = 2933.3333
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2933.3333
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2758.3333
$ws.Range("N40").Value = -4350

$ws.Range("H88").Value = 2397.7
$ws.Range("J88").Value = 2468.1428
$ws.Range("L88").Value = 2468.1428
$ws.Range("N88").Value = -3280.1428

$ws.Range("H91").Value = 2397.7
$ws.Range("J91").Value = 2468.1428
$ws.Range("L91").Value = 2468.1428
$ws.Range("N91").Value = -5276.1428

$ws.Range("H100").Value = 2473.3125
$ws.Range("I100").Value = 1324.8182
$ws.Range("K100").Value = 1324.8182
$ws.Range("M100").Value = -783.8181999999999

$ws.Range("H103").Value = 1521.8
$ws.Range("I103").Value = 1267
$ws.Range("K103").Value = 3801
$ws.Range("M103").Value = -3215

$ws.Range("H113").Value = 5659.625
$ws.Range("I113").Value = 7967.2905
$ws.Range("K113").Value = 7967.2905
$ws.Range("M113").Value = -4713.2905

$ws.Range("H129").Value = 805.5714
$ws.Range("J129").Value = 1999
$ws.Range("L129").Value = 5997
$ws.Range("N129").Value = -15997

$ws.Range("H132").Value = 1453.8
$ws.Range("I132").Value = 1416.5532
$ws.Range("K132").Value = 4249.6596
$ws.Range("M132").Value = -1719.6596

$ws.Range("H135").Value = 2446.3333
$ws.Range("I135").Value = 1135.909
$ws.Range("K135").Value = 10223.181
$ws.Range("M135").Value = -7688.181

$ws.Range("H138").Value = 2462.4749
$ws.Range("J138").Value = 2698.5422
$ws.Range("L138").Value = 8095.6266
$ws.Range("N138").Value = -18375.6266

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35950856
$ws.Range("I32").Value = 37054252
$ws.Range("K32").Value = 37054252
$ws.Range("M32").Value = -37053965

$ws.Range("H61").Value = 3508
$ws.Range("I61").Value = 3233.5715
$ws.Range("K61").Value = 3233.5715
$ws.Range("M61").Value = -3021.5715

$ws.Range("H136").Value = 3508
$ws.Range("I136").Value = 3233.5715
$ws.Range("K136").Value = 9700.7145
$ws.Range("M136").Value = -7150.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2626.6
$ws.Range("I86").Value = 3124
$ws.Range("J86").Value = 1466
$ws.Range("K86").Value = 3124
$ws.Range("L86").Value = 1466
$ws.Range("M86").Value = -2001
$ws.Range("N86").Value = -3712

$ws.Range("H89").Value = 2626.6
$ws.Range("I89").Value = 3124
$ws.Range("J89").Value = 1466
$ws.Range("K89").Value = 15620
$ws.Range("L89").Value = 7330
$ws.Range("M89").Value = -10004
$ws.Range("N89").Value = -18562

$ws.Range("H105").Value = 1997.2174
$ws.Range("I105").Value = 1390.0625
$ws.Range("K105").Value = 1390.0625
$ws.Range("M105").Value = 356.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 986.19354
$ws.Range("I107").Value = 1108.7646
$ws.Range("J107").Value = 837.3570999999999
$ws.Range("K107").Value = 3326.2938
$ws.Range("L107").Value = 2512.0713
$ws.Range("M107").Value = -1406.2938
$ws.Range("N107").Value = -6352.0713

$ws.Range("H132").Value = 386176.22
$ws.Range("I132").Value = 1110.1765
$ws.Range("J132").Value = 1113523.2
$ws.Range("K132").Value = 9991.5885
$ws.Range("L132").Value = 10021708.8
$ws.Range("M132").Value = -7461.5885
$ws.Range("N132").Value = -10026768.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 117000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 117000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 117000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -118372

$ws.Range("H65").Value = 117000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 117000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 351000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -357864

$ws.Range("H132").Value = 2775.24
$ws.Range("I132").Value = 2447.2856
$ws.Range("J132").Value = 4497
$ws.Range("K132").Value = 7341.8568
$ws.Range("L132").Value = 13491
$ws.Range("M132").Value = -4811.8568
$ws.Range("N132").Value = -18551

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 116000
$ws.Range("J63").Value = 116000
$ws.Range("L63").Value = 116000
$ws.Range("N63").Value = -117498

$ws.Range("H66").Value = 116000
$ws.Range("J66").Value = 116000
$ws.Range("L66").Value = 348000
$ws.Range("N66").Value = -355488

$ws.Range("H68").Value = 4765.8335
$ws.Range("I68").Value = 5019
$ws.Range("K68").Value = 5019
$ws.Range("M68").Value = -4270

$ws.Range("H71").Value = 4765.8335
$ws.Range("I71").Value = 5019
$ws.Range("K71").Value = 25095
$ws.Range("M71").Value = -21351

$ws.Range("H132").Value = 4064.17
$ws.Range("I132").Value = 3972.457
$ws.Range("K132").Value = 11917.371
$ws.Range("M132").Value = -9387.370999999999

$ws.Range("H136").Value = 7776.4
$ws.Range("I136").Value = 6718.0713
$ws.Range("J136").Value = 10245.833
$ws.Range("K136").Value = 20154.2139
$ws.Range("L136").Value = 30737.499
$ws.Range("M136").Value = -17604.2139
$ws.Range("N136").Value = -35837.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 54131
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 54131
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 54131
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -57625

$ws.Range("H132").Value = 1656.8857
$ws.Range("I132").Value = 1475.6897
$ws.Range("J132").Value = 2532.6667
$ws.Range("K132").Value = 4427.0691
$ws.Range("L132").Value = 7598.000100000001
$ws.Range("M132").Value = -1897.0691
$ws.Range("N132").Value = -12658.0001

$ws.Range("H136").Value = 8041469.5
$ws.Range("I136").Value = 1091.3158
$ws.Range("K136").Value = 3273.9474
$ws.Range("M136").Value = -723.9474
